# Apply the "fix p3 and p4" update: new D-column values for rows 2-6
# and move the active selection from D6 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 172.52
$ws.Range("D3").Value = 97.7
$ws.Range("D4").Value = 141.12
$ws.Range("D5").Value = 132.13999999999999
$ws.Range("D6").Value = 101.5

$ws.Range("D7").Select()
